$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Row 32 new values (previously row 33's species data, with its own new B/Taxonsorteringsordning value)
$ws.Range("A32").Value = 112213272
$ws.Range("B32").Value = 89553
$ws.Range("D32").Value = "NT"
$ws.Range("E32").Value = 1202
$ws.Range("F32").Value = "Ullticka"
$ws.Range("G32").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H32").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q32").Value = 515738
$ws.Range("R32").Value = 6704726

# Row 33 new values (previously row 32's species data, with its own new B/Taxonsorteringsordning value)
$ws.Range("A33").Value = 112213305
$ws.Range("B33").Value = 89517
$ws.Range("D33").Value = "LC"
$ws.Range("E33").Value = 5447
$ws.Range("F33").Value = "Vedticka"
$ws.Range("G33").Value = "Fuscoporia viticola"
$ws.Range("H33").Value = "(Schwein.) Murrill"
$ws.Range("Q33").Value = 515748
$ws.Range("R33").Value = 6704727

# Row 34: only Taxonsorteringsordning (column B) changes
$ws.Range("B34").Value = 90814
